$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row coloring (commit: "added colors to rows") ---
# Orange (FFCC66) rows: 7, 8, 12, 13, 15
$orange = 6737151   # RGB(255,204,102) -> FFCC66
$red    = 6184671   # RGB(223,94,94)   -> DF5E5E

$orangeRows = @(7, 8, 12, 13, 15)
foreach ($r in $orangeRows) {
    $rng = $ws.Range("A" + $r + ":J" + $r)
    $rng.Interior.Color = $orange
}

$redRows = @(14)
foreach ($r in $redRows) {
    $rng = $ws.Range("A" + $r + ":J" + $r)
    $rng.Interior.Color = $red
}

# --- Value fixes ---
# I14: 0 -> 1
$ws.Cells.Item(14, 9).Value = 1

# B19: shared-string " " -> boolean FALSE
$ws.Cells.Item(19, 2).Value = $false

# --- Formula fixes: drop the bogus extra ",1" FLOOR argument ---
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
